$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #5 closed
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.36   # Total P&L %
$summary.Range("B6").Value = 5      # Total Trades
$summary.Range("B9").Value = 60     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet: MarketMaking row picks up the new trade
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 5       # Trades
$status.Range("G4").Value = 60      # Win Rate %

# ---------------------------------------------------------------
# Helper that appends the closed trade #5 row to a trades sheet
# ---------------------------------------------------------------
function Add-TradeFiveRow($ws) {
    $ws.Cells.Item(6, 1).Value = 5
    # Force column B to remain plain text so "2026-02-17" isn't
    # auto-converted into a date serial value.
    $ws.Cells.Item(6, 2).NumberFormat = "@"
    $ws.Cells.Item(6, 2).Value = "2026-02-17"
    $ws.Cells.Item(6, 3).Value = "19:43:00"
    $ws.Cells.Item(6, 4).Value = "MarketMaking"
    $ws.Cells.Item(6, 5).Value = "DOWN"
    $ws.Cells.Item(6, 6).Value = 0.4
    $ws.Cells.Item(6, 7).Value = 0.4
    $ws.Cells.Item(6, 8).Value = "CLOSED"
    $ws.Cells.Item(6, 9).Value = 0
    $ws.Cells.Item(6, 10).Value = 0
    $ws.Cells.Item(6, 11).Value = 100.09
    $ws.Cells.Item(6, 12).Value = 0
    $ws.Cells.Item(6, 13).Value = 0
    $ws.Cells.Item(6, 14).Value = 0.6
    $ws.Cells.Item(6, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(6, 16).Value = "early_exit"
    $ws.Cells.Item(6, 17).Value = 0.13
}

# ---------------------------------------------------------------
# All Trades sheet: append the new closed trade as row 6
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeFiveRow $allTrades

# ---------------------------------------------------------------
# MarketMaking sheet: append the same closed trade as row 6
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeFiveRow $marketMaking

Write-Output "Applied trade #5 close updates"
